$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.91966724395752
$ws.Range("B1").Value = 5.435782432556152
$ws.Range("C1").Value = 4.4125075340271
$ws.Range("D1").Value = 5.153647422790527
$ws.Range("E1").Value = 4.764707088470459
